$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row containing "GNDI3" (row 10) - ticker no longer exists
$ws.Rows.Item(10).Delete()

# Fix ticker typo: "PASS3" doesn't exist, correct ticker is "PSSA3"
$found = $ws.Range("A1:A14").Find("PASS3")
if ($found -ne $null) {
    $found.Value = "PSSA3"
}

# Update the active cell selection to match the post-edit state (whole row 10 selected)
$ws.Range("A10:XFD10").Select()
